# Stabilize course search and reorganize test suites
$wb = $excel.ActiveWorkbook

# --- 1. "courses" sheet: remove the prereq_tier column (P), shifting
#        bucket1..bucket4 (Q..T) left by one column (P..S) ---
$courses = $wb.Worksheets.Item("courses")
$courses.Range("P:P").Delete()

# --- 2. "courses" sheet: populate the new bucket-count column (O) for
#        rows 46-65, which previously had no values beyond column L ---
$oValues = @{
    46 = 1
    47 = 2
    48 = 4
    49 = 2
    50 = 2
    51 = 3
    52 = 0
    53 = 2
    54 = 2
    55 = 3
    56 = 3
    57 = 3
    58 = 1
    59 = 1
    60 = 2
    61 = 2
    62 = 1
    63 = 2
    64 = 0
    65 = 2
}
foreach ($row in $oValues.Keys) {
    $courses.Cells.Item($row, 15).Value = $oValues[$row]
}

# --- 3. "tracks" sheet: clear the stray empty-string cell at E2 ---
$tracks = $wb.Worksheets.Item("tracks")
$tracks.Range("E2").ClearContents()
